# 2020.06.02 - First commit of full code reconstruction
# Adds a new LCD message row (row 7): "FUNZIONE NON IMPLEM." spread across
# the A:H / J:L / N:T column groups (mirrors the layout of rows 2-6), and
# moves the active-cell selection down to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate row 7 with the new message text. The fill order below reproduces
# the shared-string table append order seen in the saved workbook (new
# unique single-letter strings are interned in the order: O, U, N, Z, I, P, L).
$ws.Range("F7").Value = "O"
$ws.Range("B7").Value = "U"
$ws.Range("C7").Value = "N"
$ws.Range("D7").Value = "Z"
$ws.Range("E7").Value = "I"
$ws.Range("A7").Value = "F"
$ws.Range("G7").Value = "N"
$ws.Range("H7").Value = "E"
$ws.Range("J7").Value = "N"
$ws.Range("K7").Value = "O"
$ws.Range("L7").Value = "N"
$ws.Range("N7").Value = "I"
$ws.Range("O7").Value = "M"
$ws.Range("P7").Value = "P"
$ws.Range("Q7").Value = "L"
$ws.Range("R7").Value = "E"
$ws.Range("S7").Value = "M"
$ws.Range("T7").Value = "."

# Move the sheet's active cell / selection from A7 down to A8.
$ws.Range("A8").Select()
